$wb = $excel.ActiveWorkbook

# New "Kadastro" record to append as row 2 on the "Kayitlar" sheet and on the
# matching district sheet ("Merkez İlçe") — both currently hold only the
# header row.
$sheetNames = @("Kayitlar", "Merkez İlçe")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Leading "'" forces text storage (matches the workbook's existing
    # convention of storing every column, including numeric-looking ones,
    # as text) instead of letting Excel infer a number/date type.
    $ws.Cells.Item(2, 1).Value = "'1"
    $ws.Cells.Item(2, 2).Value = "'2025-08-05"
    $ws.Cells.Item(2, 3).Value = "Merkez İlçe"
    $ws.Cells.Item(2, 4).Value = "'1"
    $ws.Cells.Item(2, 5).Value = "'5"
    $ws.Cells.Item(2, 6).Value = "DÜZELTME"
    $ws.Cells.Item(2, 7).Value = "EMİNE ALANLI KIRCILI (K.Mühendisi), BARIŞ YAYLAGÜL (Tekniker)"

    # Drop the implicit "quote prefix" formatting that the apostrophe trick
    # applies, so the new row keeps the sheet's plain default style.
    $ws.Range("A2:G2").Style = "Normal"
}
